$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.812.59'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.350.73'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.04%  '
$ws.Range("D9").Value = '2.349.34'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("E11").Value = '  +2.02%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").Value = '2.775.83'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '60.652.56'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '2.347.97'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("E25").Value = '  -8.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.14%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '496.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.35%  '
$ws.Range("D31").Value = '0.0₃0866'
$ws.Range("E31").Value = '  -6.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.147'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.75%  '
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("E38").Value = '  +3.13%  '
$ws.Range("E39").Value = '  +5.95%  '
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '144.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.08%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.51%  '
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0516'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.568'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.18%  '
